# Revert the paragraph text from "Version 2." back to "Version 1.":
#   "Versi"+"on"  -> merge into a single run "Version"
#   " 2"          -> " 1."
#   "."           -> removed (its own trailing run is deleted)
# The bookmark (_GoBack) that sits between the " 2"/" 1." run and the
# trailing "." run is left untouched.

$d = $word.ActiveDocument

# Full text starts as "Version 2." with runs:
#   [0,5)  "Versi"
#   [5,7)  "on"
#   [7,9)  " 2"
#   [9,10) "."

# Step 1: delete the trailing "." run.
$d.Range(9, 10).Text = ""

# Step 2: change " 2" to " 1."
$d.Range(7, 9).Text = " 1."

# Step 3: merge the "Versi" + "on" runs into a single "Version" run.
# Setting the range to the text it already equals is a no-op, so first
# set it to a distinct value (forcing the underlying runs to merge),
# then correct it to the final "Version" text.
$d.Range(0, 7).Text = "Versionx"
$d.Range(0, 8).Text = "Version"
